$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = 8.6913127284582892
$ws.Range("C5").Value = 0.12517291953906468
$ws.Range("D5").Value = 0.83399408697976829
$ws.Range("B7").Value = 133.7687712470339
$ws.Range("C7").Value = 1.7543313225850714
$ws.Range("D7").Value = 15.253770345452764
$ws.Range("B8").Value = 16.696310127304844
$ws.Range("C8").Value = 0.21469404025735703
$ws.Range("D8").Value = 2.2402968415264843
$ws.Range("B11").Value = 26.553643778690738
$ws.Range("C11").Value = 0.40997341068083221
$ws.Range("D11").Value = 2.8961275578557197
$ws.Range("B12").Value = 79.540450382533777
$ws.Range("C12").Value = 0.98099614412832647
$ws.Range("D12").Value = 13.873181445464581
$ws.Range("B14").Value = 47.654684937898843
$ws.Range("C14").Value = 0.51526071333377077
$ws.Range("D14").Value = 9.1892542195054574
$ws.Range("B17").Value = 11.135655219222187
$ws.Range("C17").Value = 0.10261313101997194
$ws.Range("D17").Value = 2.1552679183928332
$ws.Range("B18").Value = 16.440372176267971
$ws.Range("C18").Value = 0.11623681051361495
$ws.Range("D18").Value = 5.7578986024222276
$ws.Range("B19").Value = 19.768970999658972
$ws.Range("C19").Value = 0.18916613913718447
$ws.Range("D19").Value = 5.6459004154887724
$ws.Range("B20").Value = 189.72085903267197
$ws.Range("C20").Value = 2.3167091542985516
$ws.Range("D20").Value = 34.970047880861472
$ws.Range("B21").Value = 33.406829799898354
$ws.Range("C21").Value = 0.37447281636917745
$ws.Range("D21").Value = 5.5676362236405073
$ws.Range("B22").Value = 13.571450136632571
$ws.Range("C22").Value = 0.13045316333218349
$ws.Range("D22").Value = 2.9297302287598943
$ws.Range("B23").Value = 56.120297353338721
$ws.Range("C23").Value = 0.83006082353211819
$ws.Range("D23").Value = 10.561868556765324
$ws.Range("B24").Value = 10.768622743778128
$ws.Range("C24").Value = 0.16121572456643563
$ws.Range("D24").Value = 0.80816944073490404
$ws.Range("B25").Value = 29.7244665396729
$ws.Range("C25").Value = 0.10673022586562905
$ws.Range("D25").Value = 11.323954593261602
$ws.Range("B26").Value = 53.207981146109184
$ws.Range("C26").Value = 0.90809870042147145
$ws.Range("D26").Value = 8.7311756167638421
$ws.Range("B27").Value = 75.728951208085405
$ws.Range("C27").Value = 0.90367705953825728
$ws.Range("D27").Value = 9.8404284423735149
$ws.Range("B28").Value = 47.62852356759921
$ws.Range("C28").Value = 0.4332630045433189
$ws.Range("D28").Value = 9.6523100024630093
